# Fix miss latency calculation: update column I (CPI using Miss Ticks) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 1.04807073996
$ws.Range("I3").Value = 1.04813537996
$ws.Range("I4").Value = 1.04806716996
$ws.Range("I5").Value = 1.04808419996
$ws.Range("I6").Value = 1.04808310996
$ws.Range("I7").Value = 1.04809203996
$ws.Range("I8").Value = 1.04807861996
$ws.Range("I9").Value = 1.04810255996
$ws.Range("I10").Value = 1.50774336986
$ws.Range("I11").Value = 1.19823244964
$ws.Range("I12").Value = 1.09302550996
$ws.Range("I13").Value = 1.04808419996
$ws.Range("I14").Value = 1.03522892996
$ws.Range("I15").Value = 1.02029927996
$ws.Range("I16").Value = 1.01672505996
$ws.Range("I17").Value = 1.01074883996
$ws.Range("I18").Value = 1.04819701996
$ws.Range("I19").Value = 1.04808419996
$ws.Range("I20").Value = 1.04808419996
$ws.Range("I21").Value = 1.04808419996
$ws.Range("I22").Value = 1.04808419996
$ws.Range("I23").Value = 1.04808419996
$ws.Range("I24").Value = 1.04808419996
$ws.Range("I25").Value = 1.04808419996
$ws.Range("I26").Value = 1.17459711996
$ws.Range("I27").Value = 1.10659202996
$ws.Range("I28").Value = 1.06280415996
$ws.Range("I29").Value = 1.04808419996
$ws.Range("I30").Value = 1.05173942996
$ws.Range("I31").Value = 1.07410398996
$ws.Range("I32").Value = 1.15543669996
$ws.Range("I33").Value = 1.04805311996
$ws.Range("I34").Value = 1.04808419996
$ws.Range("I35").Value = 1.04809213996
$ws.Range("I36").Value = 1.04803992996
$ws.Range("I37").Value = 1.30060880996
$ws.Range("I38").Value = 1.04808419996
$ws.Range("I39").Value = 1.02370865996
$ws.Range("I40").Value = 1.02273897996
$ws.Range("I41").Value = 1.04808419996
$ws.Range("I42").Value = 1.04808419996
$ws.Range("I43").Value = 1.04808419996
$ws.Range("I44").Value = 1.04808419996
$ws.Range("I45").Value = 11.65405681996
$ws.Range("I46").Value = 11.66127735996
$ws.Range("I47").Value = 11.65870669996
$ws.Range("I48").Value = 11.65366083996
$ws.Range("I49").Value = 11.65132341996
$ws.Range("I50").Value = 11.64757588996
$ws.Range("I51").Value = 11.66132641996
$ws.Range("I52").Value = 11.66179159996
$ws.Range("I53").Value = 11.70041760996
$ws.Range("I54").Value = 11.66172645996
$ws.Range("I55").Value = 11.66527231996
$ws.Range("I56").Value = 11.65366083996
$ws.Range("I57").Value = 11.65074453996
$ws.Range("I58").Value = 11.67263649996
$ws.Range("I59").Value = 11.67263585996
$ws.Range("I60").Value = 11.67263562996
$ws.Range("I61").Value = 11.67408785996
$ws.Range("I62").Value = 11.64913541996
$ws.Range("I63").Value = 11.67744037996
$ws.Range("I64").Value = 11.65366083996
$ws.Range("I65").Value = 11.63772217996
$ws.Range("I66").Value = 11.60592929996
$ws.Range("I67").Value = 11.54750929996
$ws.Range("I68").Value = 9.76632795996
$ws.Range("I69").Value = 33.16575241996
$ws.Range("I70").Value = 41.16394051996
$ws.Range("I71").Value = 21.53758099996
$ws.Range("I72").Value = 11.65366083996
$ws.Range("I73").Value = 7.030991089960001
$ws.Range("I74").Value = 4.98849250996
$ws.Range("I75").Value = 3.29345817996
$ws.Range("I76").Value = 11.64876733996
$ws.Range("I77").Value = 11.65366083996
$ws.Range("I78").Value = 11.65889467996
$ws.Range("I79").Value = 11.65931334996
$ws.Range("I80").Value = 11.66430012996
$ws.Range("I81").Value = 11.65366083996
$ws.Range("I82").Value = 11.65199244996
$ws.Range("I83").Value = 11.65199651996
$ws.Range("I84").Value = 11.64132725996
$ws.Range("I85").Value = 11.65366083996
$ws.Range("I86").Value = 11.67336391996
$ws.Range("I87").Value = 11.65775287996
